$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.232.14'
$ws.Range('E2').Value = '  -0.15%  '
$ws.Range('D3').Value = '1.894.75'
$ws.Range('E3').Value = '  +1.58%  '
$ws.Range('E4').Value = '  -0.28%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '245.49'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  +2.53%  '
$ws.Range('E6').Value = '  +5.84%  '
$ws.Range('E7').Value = '  -0.17%  '
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '41.10'
$c.Style = 'Normal'
$ws.Range('E8').Value = '  -3.38%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.348'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  +5.23%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '52.68'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  +12.20%  '
$ws.Range('E11').Value = '  +2.70%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '0.0995'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  +0.48%  '
$ws.Range('D13').Value = '2.170.21'
$ws.Range('E13').Value = '  +1.63%  '
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '12.08'
$c.Style = 'Normal'
$ws.Range('E14').Value = '  +4.50%  '
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '0.695'
$c.Style = 'Normal'
$ws.Range('E15').Value = '  +2.08%  '
$ws.Range('D16').Value = '1.887.58'
$ws.Range('E16').Value = '  +1.15%  '
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '4.83'
$c.Style = 'Normal'
$ws.Range('E17').Value = '  +1.99%  '
$ws.Range('D18').Value = '35.229.35'
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '72.31'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  +3.02%  '
$ws.Range('D20').Value = '0.0₃0818'
$ws.Range('E20').Value = '  +2.56%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '239.92'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  -0.58%  '
$ws.Range('E22').Value = '  +1.08%  '
$ws.Range('E23').Value = '  +1.49%  '
$ws.Range('E24').Value = '  -0.15%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '2.30'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  +1.43%  '
$ws.Range('E26').Value = '  +21.59%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '170.27'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  +0.32%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '8.41'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  +3.48%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '18.29'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  +3.04%  '
$ws.Range('E30').Value = '  +1.27%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '4.11'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  +2.15%  '
$ws.Range('E32').Value = '  -0.15%  '
$ws.Range('E33').Value = '  -0.24%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '0.931'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  +13.72%  '
$ws.Range('E35').Value = '  +0.82%  '
$ws.Range('E36').Value = '  -4.91%  '
$ws.Range('E37').Value = '  -3.30%  '
$ws.Range('E38').Value = '  +1.16%  '
$ws.Range('E39').Value = '  -1.75%  '
$ws.Range('E40').Value = '  +2.20%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '15.96'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  +4.45%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '0.0626'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  +3.94%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '89.15'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  -1.44%  '
$ws.Range('D44').Value = '1.336.31'
$ws.Range('E44').Value = '  -0.74%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '2.37'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  +0.81%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '47.87'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  +37.21%  '
$ws.Range('E47').Value = '  -0.55%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '2.76'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  +0.79%  '
$ws.Range('E49').Value = '  -2.11%  '
$ws.Range('B50').Value = 'RocketPoolETH'
$ws.Range('C50').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D50').Value = '2.080.60'
$ws.Range('E50').Value = '  +1.42%  '
$ws.Range('B51').Value = 'Gas'
$ws.Range('C51').Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '11.54'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  -6.95%  '
